{"js": "// Update the addressee name, street, postal code/city, and the letter date\n// in the reminder (\"Mahnung\") template.\nconst replacements = [\n  [\"David Achermann\", \"Reto Affolter\"],\n  [\"S\u00fcdringstrasse 73\", \"H\u00fcseliring 12\"],\n  [\"4563 Gerlafingen\", \"4565 Recherswil\"],\n  [\"18.01.2025\", \"03.06.2025\"]\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the addressee name, street, postal code/city, and the letter date\n# in the reminder (\"Mahnung\") template.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"David Achermann\", \"Reto Affolter\"),\n    @(\"S\u00fcdringstrasse 73\", \"H\u00fcseliring 12\"),\n    @(\"4563 Gerlafingen\", \"4565 Recherswil\"),\n    @(\"18.01.2025\", \"03.06.2025\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
